$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Unprotect()

$ws.Range("A18").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-26 for illustrative purposes only and are subject to change."

$ws.Range("D2").Value = 0.0565299256305414
$ws.Range("E2").Value = 0.003602105846494785

$ws.Range("D3").Value = 0.02383412249586383

$ws.Range("D4").Value = 0.03119601361427737
$ws.Range("E4").Value = 0.00734324985878354

$ws.Range("D5").Value = 0.03207970160845573
$ws.Range("E5").Value = 0.006957866254348577

$ws.Range("D6").Value = 0.03593229409710778
$ws.Range("E6").Value = 0.005231037489101809

$ws.Range("D7").Value = 0.01867166161406852
$ws.Range("E7").Value = 0.002446782481037468

$ws.Range("D8").Value = 0.004417134671048122
$ws.Range("E8").Value = 0.01196808510638303

$ws.Range("D9").Value = 0.006864571877967517
$ws.Range("E9").Value = 0.003612854154782186

$ws.Range("D10").Value = 0.07393218314662106
$ws.Range("E10").Value = -0.001588983050847426

$ws.Range("D11").Value = 0.07404966013255318
$ws.Range("E11").Value = -0.001586462189317794

$ws.Range("D12").Value = 0.1456296929608443
$ws.Range("E12").Value = -0.001720923562311794

$ws.Range("D13").Value = 0.3819137653658266
$ws.Range("E13").Value = -0.00113299633955033

$ws.Range("D14").Value = 0.1149492727848246
$ws.Range("E14").Value = -0.004833671157744024

$ws.Range("E15").Value = -0.0005066782403252423

$ws.Protect()
